# Add a new "2020-12-27" attendance sheet at the end of the workbook,
# matching the layout/styling already used by the other daily-attendance
# sheets (e.g. "2020-12-21").

$wb = $excel.ActiveWorkbook

# Remember the sheet that is active before we start, so we can restore the
# original selection/active-tab once we're done adding the new sheet.
$origActive = $wb.ActiveSheet

# Insert the new worksheet after the current last sheet, so it lands at the
# very end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2020-12-27"

# -- Header row ------------------------------------------------------------
$ws.Range("A1").Value = "Sr. No"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Address"
$ws.Range("D1").Value = "Job"
$ws.Range("E1").Value = "Time-Stamp"
$ws.Range("F1").Value = "SpO2_value"
$ws.Range("G1").Value = "Heart-rate"
$ws.Range("H1").Value = "Compensated"
$ws.Range("I1").Value = "Ambient"

# -- Data row ----------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "sachin"
$ws.Range("C2").Value = "301/Sanskruti-1,Andheri, Mumbai"
$ws.Range("D2").Value = "Software Engineer"
$ws.Range("E2").Value = "10:45:04"
$ws.Range("F2").Value = 95.83510336095262
$ws.Range("G2").Value = 67.90663936802549
$ws.Range("H2").Value = "NA"
$ws.Range("I2").Value = "NA"

# -- Formatting: header row (A1:I1) + the "Sr. No" cell (A2) are bold,
#    centered, top-aligned, boxed with a thin border -- same look as the
#    sibling sheets.
$headerRange = $ws.Range("A1:I1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous

$srNoCell = $ws.Range("A2")
$srNoCell.Font.Bold = $true
$srNoCell.HorizontalAlignment = -4108      # xlCenter
$srNoCell.VerticalAlignment = -4160        # xlTop
$srNoCell.Borders.LineStyle = 1            # xlContinuous

# -- Page margins: 0.75"/0.75"/1"/1" (L/R/T/B), 0.5" header/footer -- same
#    as the other sheets (values are in points: 72pt = 1", 54pt = 0.75",
#    36pt = 0.5").
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Restore whichever sheet/tab was active before we added the new one, so the
# new sheet is appended without disturbing the existing active-tab state.
$origActive.Activate() | Out-Null
$origActive.Range("A1").Select() | Out-Null
